$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 988.375
$ws.Range("I32").Value = 1001
$ws.Range("J32").Value = 986.5714
$ws.Range("K32").Value = 1001
$ws.Range("L32").Value = 986.5714
$ws.Range("M32").Value = -675
$ws.Range("N32").Value = -1638.5714

$ws.Range("H43").Value = 2355.5
$ws.Range("I43").Value = 2566.6667
$ws.Range("J43").Value = 2313.2666
$ws.Range("K43").Value = 2566.6667
$ws.Range("L43").Value = 2313.2666
$ws.Range("M43").Value = -2497.6667
$ws.Range("N43").Value = -2451.2666

$ws.Range("H121").Value = 1365.7407
$ws.Range("I121").Value = 150
$ws.Range("J121").Value = 1412.5
$ws.Range("K121").Value = 450
$ws.Range("L121").Value = 4237.5
$ws.Range("M121").Value = 1297
$ws.Range("N121").Value = -7731.5

$ws.Range("H137").Value = 3474342.8
$ws.Range("I137").Value = 4631302.5
$ws.Range("J137").Value = 3463.3333
$ws.Range("K137").Value = 13893907.5
$ws.Range("L137").Value = 10389.9999
$ws.Range("M137").Value = -13891357.5
$ws.Range("N137").Value = -15489.9999

$ws.Range("H138").Value = 4246.458
$ws.Range("I138").Value = 3751
$ws.Range("J138").Value = 4393.547
$ws.Range("K138").Value = 11253
$ws.Range("L138").Value = 13180.641
$ws.Range("M138").Value = -6113
$ws.Range("N138").Value = -23460.641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17259646
$ws.Range("I32").Value = 19624372
$ws.Range("J32").Value = 30928.572
$ws.Range("K32").Value = 19624372
$ws.Range("L32").Value = 30928.572
$ws.Range("M32").Value = -19624085
$ws.Range("N32").Value = -31502.572

$ws.Range("H74").Value = 12822071
$ws.Range("I74").Value = 1143.7354
$ws.Range("J74").Value = 100004380
$ws.Range("K74").Value = 1143.7354
$ws.Range("L74").Value = 100004380
$ws.Range("M74").Value = -269.7354
$ws.Range("N74").Value = -100006128

$ws.Range("H77").Value = 12822071
$ws.Range("I77").Value = 1143.7354
$ws.Range("J77").Value = 100004380
$ws.Range("K77").Value = 5718.677
$ws.Range("L77").Value = 500021900
$ws.Range("M77").Value = -1350.677
$ws.Range("N77").Value = -500030636

$ws.Range("H93").Value = 72815.664
$ws.Range("J93").Value = 72815.664
$ws.Range("L93").Value = 72815.664
$ws.Range("N93").Value = -77807.664

$ws.Range("H95").Value = 87500
$ws.Range("J95").Value = 87500
$ws.Range("L95").Value = 87500
$ws.Range("N95").Value = -92992

$ws.Range("H132").Value = 2656920.5
$ws.Range("I132").Value = 5347.8823
$ws.Range("J132").Value = 6413315.5
$ws.Range("K132").Value = 16043.6469
$ws.Range("L132").Value = 19239946.5
$ws.Range("M132").Value = -13513.6469
$ws.Range("N132").Value = -19245006.5

$ws.Range("H137").Value = 65847.27
$ws.Range("I137").Value = 43500
$ws.Range("J137").Value = 68082
$ws.Range("K137").Value = 43500
$ws.Range("L137").Value = 68082
$ws.Range("M137").Value = -38400
$ws.Range("N137").Value = -78282

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 53608.367
$ws.Range("I107").Value = 77732.234
$ws.Range("K107").Value = 77732.234
$ws.Range("M107").Value = -75812.234

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1927.375
$ws.Range("I16").Value = 1922.3529
$ws.Range("J16").Value = 1939.5714
$ws.Range("K16").Value = 1922.3529
$ws.Range("L16").Value = 1939.5714
$ws.Range("M16").Value = -1635.3529
$ws.Range("N16").Value = -2513.5714

$ws.Range("H31").Value = 7514.3374
$ws.Range("I31").Value = 1331
$ws.Range("J31").Value = 8605.514999999999
$ws.Range("K31").Value = 1331
$ws.Range("L31").Value = 8605.514999999999
$ws.Range("M31").Value = -1036
$ws.Range("N31").Value = -9195.514999999999

$ws.Range("H34").Value = 7514.3374
$ws.Range("I34").Value = 1331
$ws.Range("J34").Value = 8605.514999999999
$ws.Range("K34").Value = 1331
$ws.Range("L34").Value = 8605.514999999999
$ws.Range("M34").Value = -1129
$ws.Range("N34").Value = -9009.514999999999

$ws.Range("H107").Value = 8929703
$ws.Range("I107").Value = 15625880
$ws.Range("J107").Value = 1466.6666
$ws.Range("K107").Value = 15625880
$ws.Range("L107").Value = 1466.6666
$ws.Range("M107").Value = -15623960
$ws.Range("N107").Value = -5306.6666

$ws.Range("H113").Value = 1927.375
$ws.Range("I113").Value = 1922.3529
$ws.Range("J113").Value = 1939.5714
$ws.Range("K113").Value = 1922.3529
$ws.Range("L113").Value = 1939.5714
$ws.Range("M113").Value = 247.6470999999999
$ws.Range("N113").Value = -6279.5714

$ws.Range("H132").Value = 37039184
$ws.Range("I132").Value = 45456464
$ws.Range("J132").Value = 23812024
$ws.Range("K132").Value = 136369392
$ws.Range("L132").Value = 71436072
$ws.Range("M132").Value = -136366862
$ws.Range("N132").Value = -71441132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 619.3484999999999
$ws.Range("I113").Value = 596.1395
$ws.Range("K113").Value = 1788.4185
$ws.Range("M113").Value = 381.5815

$ws.Range("H131").Value = 3342.6538
$ws.Range("I131").Value = 775.55554
$ws.Range("J131").Value = 3879.9534
$ws.Range("K131").Value = 2326.66662
$ws.Range("L131").Value = 11639.8602
$ws.Range("M131").Value = 2713.33338
$ws.Range("N131").Value = -21719.8602

$ws.Range("H132").Value = 2253.0483
$ws.Range("J132").Value = 2076.853
$ws.Range("L132").Value = 18691.677
$ws.Range("N132").Value = -23751.677

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9992
$ws.Range("J2").Value = 9992
$ws.Range("L2").Value = 9992
$ws.Range("N2").Value = -10216

$ws.Range("H61").Value = 2736.1875
$ws.Range("I61").Value = 1024.4546
$ws.Range("J61").Value = 6502
$ws.Range("K61").Value = 1024.4546
$ws.Range("L61").Value = 6502
$ws.Range("M61").Value = -822.4546
$ws.Range("N61").Value = -6906

$ws.Range("H113").Value = 2736.1875
$ws.Range("I113").Value = 1024.4546
$ws.Range("J113").Value = 6502
$ws.Range("K113").Value = 1024.4546
$ws.Range("L113").Value = 6502
$ws.Range("M113").Value = 1145.5454
$ws.Range("N113").Value = -10842

$ws.Range("H132").Value = 3120.5908
$ws.Range("I132").Value = 2491.5417
$ws.Range("J132").Value = 3875.45
$ws.Range("K132").Value = 7474.625100000001
$ws.Range("L132").Value = 11626.35
$ws.Range("M132").Value = -4944.625100000001
$ws.Range("N132").Value = -16686.35

$ws.Range("H137").Value = 35107.25
$ws.Range("J137").Value = 35107.25
$ws.Range("L137").Value = 35107.25
$ws.Range("N137").Value = -45307.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 899.2083
$ws.Range("I113").Value = 947.3182
$ws.Range("J113").Value = 370
$ws.Range("K113").Value = 2841.9546
$ws.Range("L113").Value = 1110
$ws.Range("M113").Value = -671.9546
$ws.Range("N113").Value = -5450

$ws.Range("H132").Value = 5119123
$ws.Range("I132").Value = 2132.389
$ws.Range("J132").Value = 13891107
$ws.Range("K132").Value = 6397.167
$ws.Range("L132").Value = 41673321
$ws.Range("M132").Value = -3867.167
$ws.Range("N132").Value = -41678381
